$p = $ppt.ActivePresentation
$m = $p.Slides.Item(1).Master
$t = $m.Theme
try {
  $t.Bogus()
} catch {
  Write-Host "ERROR:" $_.Exception.Message
}
